# Add a new "empty sheet" flag entry to the flag_map dictionary (row 56):
#   empty_sheet | Loaded template has empty sheet | Hard Stop (Empty Sheet)
# This mirrors the PKWG PCB dataset conversion dictionary update.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A56").Value = "empty_sheet"
$ws.Range("B56").Value = "Loaded template has empty sheet"
$ws.Range("C56").Value = "Hard Stop (Empty Sheet)"

# Match the author's final selection/active cell on save.
$ws.Range("C56").Select()
